$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input cells in row 7 (Test_2) with new measured values
$ws.Range("E7").Value = 58.51999999999996
$ws.Range("G7").Value = 5.319999999999931
$ws.Range("I7").Value = 4.4499999999999496

$excel.CalculateFull()
